$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.403.88'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '3.136.66'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = "'242.60"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.84%  '
$ws.Range('D6').Value = "'627.40"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  +9.58%  '
$ws.Range('D8').Value = "'0.372"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.31%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '3.134.53'
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').Value = "'0.766"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.53%  '
$ws.Range('E12').Value = '  +4.62%  '
$ws.Range('D13').Value = "'0.0000253"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.89%  '
$ws.Range('D14').Value = "'35.88"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = "'5.51"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.96%  '
$ws.Range('D16').Value = '90.799.20'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '3.725.25'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = '3.114.87'
$ws.Range('E18').Value = '  -2.85%  '
$ws.Range('D19').Value = "'3.78"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.75%  '
$ws.Range('D20').Value = "'14.68"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('D21').Value = "'0.0000213"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = "'5.87"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.79%  '
$ws.Range('D23').Value = "'451.97"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').Value = "'9.17"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D25').Value = "'5.97"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.49%  '
$ws.Range('D26').Value = "'93.38"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.90%  '
$ws.Range('D27').Value = "'11.99"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.52%  '
$ws.Range('D28').Value = '3.297.80'
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = "'0.179"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.76%  '
$ws.Range('D31').Value = "'0.125"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +43.62%  '
$ws.Range('D32').Value = "'0.230"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +14.15%  '
$ws.Range('D33').Value = "'9.12"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.34%  '
$ws.Range('E34').Value = '  +25.51%  '
$ws.Range('D35').Value = "'0.163"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.69%  '
$ws.Range('D36').Value = "'27.02"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.31%  '
$ws.Range('D37').Value = "'7.66"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.07%  '
$ws.Range('D38').Value = "'4.20"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +24.77%  '
$ws.Range('D39').Value = "'503.10"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('D40').Value = "'1.93"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = "'3.62"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.48%  '
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('D43').Value = "'0.426"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').Value = "'22.17"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('D47').Value = "'157.31"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.33%  '
$ws.Range('D48').Value = "'0.700"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('D49').Value = "'4.58"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = "'45.14"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').Value = "'1.35"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.37%  '
